$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 940
$ws.Range("F5").Value = 883
$ws.Range("F6").Value = 476
$ws.Range("F9").Value = 38818
$ws.Range("G9").Value = "暂时售罄"
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 149
$ws.Range("F12").Value = 518
$ws.Range("F15").Value = 0
$ws.Range("F17").Value = 186
$ws.Range("F18").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 197
$ws.Range("F23").Value = 1007
$ws.Range("F25").Value = 531
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 550
$ws.Range("F30").Value = 14
$ws.Range("F33").Value = 807
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 217
$ws.Range("F39").Value = 41
$ws.Range("F40").Value = 974
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("F44").Value = 0
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 186
$ws.Range("F5").Value = 4356
$ws.Range("F8").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F15").Value = 67
$ws.Range("F17").Value = 2
$ws.Range("F18").Value = 4354
$ws.Range("F21").Value = 5
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 258
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 258
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 940
$ws.Range("F7").Value = 91
$ws.Range("F8").Value = 883
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = 0
$ws.Range("G11").Value = "暂时售罄"
$ws.Range("F12").Value = 186
$ws.Range("F13").Value = 322
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 8143
$ws.Range("F18").Value = 149
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("F23").Value = 99
$ws.Range("F24").Value = 186
$ws.Range("F25").Value = 619
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 44
$ws.Range("F29").Value = 1007
$ws.Range("F30").Value = 297
$ws.Range("F32").Value = 0
$ws.Range("F34").Value = 550
$ws.Range("F35").Value = 32
$ws.Range("F36").Value = 0
$ws.Range("F38").Value = 119
$ws.Range("F44").Value = 183
$ws.Range("F46").Value = 323
$ws.Range("F49").Value = 0
